$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 123.24
$ws.Range("I33").Value = 123.24
$ws.Range("K33").Value = 123.24
$ws.Range("M33").Value = 105.76

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 9975
$ws.Range("J3").Value = 9975
$ws.Range("L3").Value = 9975
$ws.Range("N3").Value = -10205
$ws.Range("H11").Value = 14999.5
$ws.Range("J11").Value = 20000
$ws.Range("L11").Value = 20000
$ws.Range("N11").Value = -20288
$ws.Range("H74").Value = 3143.889
$ws.Range("I74").Value = 2614.2856
$ws.Range("J74").Value = 4997.5
$ws.Range("K74").Value = 2614.2856
$ws.Range("L74").Value = 4997.5
$ws.Range("M74").Value = -1740.2856
$ws.Range("N74").Value = -6745.5
$ws.Range("H77").Value = 3143.889
$ws.Range("I77").Value = 2614.2856
$ws.Range("J77").Value = 4997.5
$ws.Range("K77").Value = 13071.428
$ws.Range("L77").Value = 24987.5
$ws.Range("M77").Value = -8703.428
$ws.Range("N77").Value = -33723.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H5").Value = 3450.75
$ws.Range("I5").Value = 399.5
$ws.Range("J5").Value = 6502
$ws.Range("K5").Value = 399.5
$ws.Range("L5").Value = 6502
$ws.Range("M5").Value = -286.5
$ws.Range("N5").Value = -6728
$ws.Range("H8").Value = 55376
$ws.Range("I8").Value = 105502
$ws.Range("J8").Value = 5250
$ws.Range("K8").Value = 105502
$ws.Range("L8").Value = 5250
$ws.Range("M8").Value = -105362
$ws.Range("N8").Value = -5530
$ws.Range("H96").Value = 20428
$ws.Range("I96").Value = 20428
$ws.Range("K96").Value = 20428
$ws.Range("M96").Value = -17682
$ws.Range("H112").Value = 1000000
$ws.Range("J112").Value = 1000000
$ws.Range("L112").Value = 1000000
$ws.Range("N112").Value = -1002954
$ws.Range("H131").Value = 9998
$ws.Range("I131").Value = 9998
$ws.Range("K131").Value = 9998
$ws.Range("M131").Value = -4958
$ws.Range("H134").Value = 10386.8
$ws.Range("I134").Value = 10933.5
$ws.Range("J134").Value = 8200
$ws.Range("K134").Value = 32800.5
$ws.Range("L134").Value = 24600
$ws.Range("M134").Value = -30265.5
$ws.Range("N134").Value = -29670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1526.875
$ws.Range("I2").Value = 914.75
$ws.Range("J2").Value = 2139
$ws.Range("K2").Value = 914.75
$ws.Range("L2").Value = 2139
$ws.Range("M2").Value = -801.75
$ws.Range("N2").Value = -2365
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 111
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 552.5
$ws.Range("H32").Value = 594.25
$ws.Range("J32").Value = 594.25
$ws.Range("L32").Value = 1782.75
$ws.Range("N32").Value = -2348.75
$ws.Range("H34").Value = 1370
$ws.Range("J34").Value = 2009.8
$ws.Range("L34").Value = 6029.4
$ws.Range("N34").Value = -6197.4
$ws.Range("H38").Value = 117
$ws.Range("I38").Value = 120.4
$ws.Range("K38").Value = 361.2
$ws.Range("M38").Value = -14.20000000000005
$ws.Range("H39").Value = 2366.5833
$ws.Range("J39").Value = 2366.5833
$ws.Range("L39").Value = 7099.749899999999
$ws.Range("N39").Value = -7687.749899999999
$ws.Range("H68").Value = 931.3333
$ws.Range("J68").Value = 995
$ws.Range("L68").Value = 2985
$ws.Range("N68").Value = -4607
$ws.Range("H71").Value = 931.3333
$ws.Range("J71").Value = 995
$ws.Range("L71").Value = 8955
$ws.Range("N71").Value = -17067
$ws.Range("H86").Value = 1498.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1498.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 4495.5
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -6867.5
$ws.Range("H89").Value = 1498.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1498.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 13486.5
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -25342.5
$ws.Range("H92").Value = 1500
$ws.Range("J92").Value = 1500
$ws.Range("L92").Value = 4500
$ws.Range("N92").Value = -6996
$ws.Range("H94").Value = 3766.6667
$ws.Range("H128").Value = 307284.56
$ws.Range("I128").Value = 307284.56
$ws.Range("K128").Value = 921853.6799999999
$ws.Range("M128").Value = -916873.6799999999
$ws.Range("H131").Value = 1331.8
$ws.Range("I131").Value = 967.8
$ws.Range("J131").Value = 1695.8
$ws.Range("K131").Value = 2903.4
$ws.Range("L131").Value = 5087.4
$ws.Range("M131").Value = 2136.6
$ws.Range("N131").Value = -15167.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 900
$ws.Range("I80").Value = 900
$ws.Range("K80").Value = 900
$ws.Range("M80").Value = 98
$ws.Range("H83").Value = 900
$ws.Range("I83").Value = 900
$ws.Range("K83").Value = 4500
$ws.Range("M83").Value = 492
$ws.Range("H132").Value = 6000
$ws.Range("I132").Value = 6000
$ws.Range("K132").Value = 18000
$ws.Range("M132").Value = -15470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 886.4706
$ws.Range("I22").Value = 762
$ws.Range("K22").Value = 762
$ws.Range("M22").Value = -467
$ws.Range("H27").Value = 886.4706
$ws.Range("I27").Value = 762
$ws.Range("K27").Value = 762
$ws.Range("M27").Value = -655
$ws.Range("H132").Value = 2097.8
$ws.Range("I132").Value = 2097.8
$ws.Range("K132").Value = 6293.400000000001
$ws.Range("M132").Value = -3763.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 84200.2
$ws.Range("I2").Value = 8001.5
$ws.Range("J2").Value = 134999.33
$ws.Range("K2").Value = 8001.5
$ws.Range("L2").Value = 134999.33
$ws.Range("M2").Value = -7889.5
$ws.Range("N2").Value = -135223.33
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H68").Value = 52500
$ws.Range("J68").Value = 52500
$ws.Range("L68").Value = 52500
$ws.Range("N68").Value = -54122
$ws.Range("H71").Value = 52500
$ws.Range("J71").Value = 52500
$ws.Range("L71").Value = 157500
$ws.Range("N71").Value = -165612
$ws.Range("H132").Value = 2486.889
$ws.Range("I132").Value = 1912.1428
$ws.Range("J132").Value = 4498.5
$ws.Range("K132").Value = 5736.428400000001
$ws.Range("L132").Value = 13495.5
$ws.Range("M132").Value = -3206.428400000001
$ws.Range("N132").Value = -18555.5
